$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix on employee in bulk: update the sample employee row with
# correct test data (firstName/lastName/email).
$ws.Range("B2").Value = "leon"
$ws.Range("A2").Value = "kibdne"
$ws.Range("C2").Value = "kinde@gmail.com"

# Move active selection to D6 (matches the saved cursor position)
$ws.Range("D6").Select()
